$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) on input sheet
$wsInput.Range("B1").Value = "2455-RBI-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"

# Update short name (B2) on input sheet - now stored as text instead of a number
$wsInput.Range("B2").Value = "245d"

# Mirror the product name onto the output sheet
$wsOutput.Range("B1").Value = "2455-RBI-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"

# Update the active selection on the input sheet
$wsInput.Activate()
$wsInput.Range("B1").Select()
